# Generate Report for Handback
# Adds a new handed-back file (2f7db598-7b81-4391-a6ab-0a7ed8fae673.md) as row 3
# on each of the three report sheets, and refreshes the existing row (2)
# which is renamed from 96555b91-741f-48b4-9887-5c2f343ac0d9.md to
# 27f5425e-7377-4959-9110-1f54699a9831.md with newer timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Overview sheet - update existing row 2, then add row 3
# ---------------------------------------------------------------------

$wsOverview.Range("A2").Value = "27f5425e-7377-4959-9110-1f54699a9831.md"
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-25 09:06:15"

$wsOverview.Range("A3").Value = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-25 09:06:15"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd835ac9d00d492a75c84b10326b7ffcac005ab4/e2e/27f5425e-7377-4959-9110-1f54699a9831.md", "", "", "e2e\27f5425e-7377-4959-9110-1f54699a9831.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd835ac9d00d492a75c84b10326b7ffcac005ab4/e2e/2f7db598-7b81-4391-a6ab-0a7ed8fae673.md", "", "", "e2e\2f7db598-7b81-4391-a6ab-0a7ed8fae673.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# 2. zh-cn sheet - update existing row 2, then add row 3
# ---------------------------------------------------------------------

$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D2").Value = "e2e"
$wsZhCn.Range("E2").Value = "ht"
$wsZhCn.Range("F2").Value = "False"
$wsZhCn.Range("G2").Value = "27f5425e-7377-4959-9110-1f54699a9831.518a67538865fa12b14fb87deb7d674cd9a0135d.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-25 09:06:01"
$wsZhCn.Range("J2").Value = "27f5425e-7377-4959-9110-1f54699a9831.518a67538865fa12b14fb87deb7d674cd9a0135d.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-25 09:06:34"
$wsZhCn.Range("L2").Value = ""
$wsZhCn.Range("M2").Value = "True"
$wsZhCn.Range("N2").Value = ""
$wsZhCn.Range("O2").Value = "False"
$wsZhCn.Range("P2").Value = ""

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.99e9a99aeeb7aca0e52c690894766408abdc79e1.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 09:06:01"
$wsZhCn.Range("J3").Value = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.99e9a99aeeb7aca0e52c690894766408abdc79e1.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-25 09:06:34"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd835ac9d00d492a75c84b10326b7ffcac005ab4/e2e/27f5425e-7377-4959-9110-1f54699a9831.md", "", "", "27f5425e-7377-4959-9110-1f54699a9831.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fe3c1d42c5674356b42a82086b5290be2583c130/e2e/27f5425e-7377-4959-9110-1f54699a9831.md", "", "", "27f5425e-7377-4959-9110-1f54699a9831.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd835ac9d00d492a75c84b10326b7ffcac005ab4/e2e/2f7db598-7b81-4391-a6ab-0a7ed8fae673.md", "", "", "2f7db598-7b81-4391-a6ab-0a7ed8fae673.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fe3c1d42c5674356b42a82086b5290be2583c130/e2e/2f7db598-7b81-4391-a6ab-0a7ed8fae673.md", "", "", "2f7db598-7b81-4391-a6ab-0a7ed8fae673.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# 3. de-de sheet - update existing row 2, then add row 3
# ---------------------------------------------------------------------

$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D2").Value = "e2e"
$wsDeDe.Range("E2").Value = "ht"
$wsDeDe.Range("F2").Value = "False"
$wsDeDe.Range("G2").Value = "27f5425e-7377-4959-9110-1f54699a9831.518a67538865fa12b14fb87deb7d674cd9a0135d.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-25 09:06:15"
$wsDeDe.Range("J2").Value = "27f5425e-7377-4959-9110-1f54699a9831.518a67538865fa12b14fb87deb7d674cd9a0135d.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-25 09:06:41"
$wsDeDe.Range("L2").Value = ""
$wsDeDe.Range("M2").Value = "True"
$wsDeDe.Range("N2").Value = ""
$wsDeDe.Range("O2").Value = "False"
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.99e9a99aeeb7aca0e52c690894766408abdc79e1.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 09:06:15"
$wsDeDe.Range("J3").Value = "2f7db598-7b81-4391-a6ab-0a7ed8fae673.99e9a99aeeb7aca0e52c690894766408abdc79e1.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-25 09:06:41"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd835ac9d00d492a75c84b10326b7ffcac005ab4/e2e/27f5425e-7377-4959-9110-1f54699a9831.md", "", "", "27f5425e-7377-4959-9110-1f54699a9831.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/de0567e9a7b0220811a75e86cbe08e8f2677b8be/e2e/27f5425e-7377-4959-9110-1f54699a9831.md", "", "", "27f5425e-7377-4959-9110-1f54699a9831.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dd835ac9d00d492a75c84b10326b7ffcac005ab4/e2e/2f7db598-7b81-4391-a6ab-0a7ed8fae673.md", "", "", "2f7db598-7b81-4391-a6ab-0a7ed8fae673.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/de0567e9a7b0220811a75e86cbe08e8f2677b8be/e2e/2f7db598-7b81-4391-a6ab-0a7ed8fae673.md", "", "", "2f7db598-7b81-4391-a6ab-0a7ed8fae673.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
